# Apply latest crypto price/volume snapshot to the "cryptos" worksheet.
# (Updated cryptos list on Sat Feb  3 06:25:30 UTC 2024 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while always preserving its original
# "text" storage - the Price column holds values such as "43.158.47" or
# "1.00" which Excel would otherwise auto-convert to a number/date.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '43.128.20'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '2.319.66'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("E4").Value = '  +0.03%  '
Set-TextValue $ws.Range("D5") '303.29'
$ws.Range("E5").Value = '  +0.50%  '
Set-TextValue $ws.Range("D6") '99.76'
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +2.58%  '
Set-TextValue $ws.Range("D10") '36.23'
$ws.Range("E10").Value = '  +5.60%  '
Set-TextValue $ws.Range("D11") '0.0792'
$ws.Range("E11").Value = '  -0.78%  '
$ws.Range("E12").Value = '  -0.99%  '
Set-TextValue $ws.Range("D13") '17.74'
$ws.Range("E13").Value = '  -1.07%  '
Set-TextValue $ws.Range("D14") '6.91'
$ws.Range("E14").Value = '  +1.57%  '
$ws.Range("D15").Value = '2.680.42'
$ws.Range("E15").Value = '  +0.56%  '
$ws.Range("D16").Value = '2.311.72'
$ws.Range("E16").Value = '  -0.35%  '
Set-TextValue $ws.Range("D17") '0.796'
$ws.Range("E17").Value = '  -2.23%  '
$ws.Range("D18").Value = '43.045.52'
$ws.Range("E18").Value = '  -0.01%  '
Set-TextValue $ws.Range("D19") '13.09'
$ws.Range("E19").Value = '  +3.82%  '
Set-TextValue $ws.Range("D20") '6.20'
$ws.Range("E20").Value = '  +1.60%  '
$ws.Range("D21").Value = '0.0₃0911'
$ws.Range("E21").Value = '  +0.43%  '
Set-TextValue $ws.Range("D22") '68.29'
$ws.Range("E22").Value = '  +0.82%  '
Set-TextValue $ws.Range("D23") '240.29'
$ws.Range("E23").Value = '  +1.24%  '
Set-TextValue $ws.Range("D24") '2.16'
$ws.Range("E24").Value = '  -2.35%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  -0.12%  '
Set-TextValue $ws.Range("D27") '25.56'
$ws.Range("E27").Value = '  +3.26%  '
Set-TextValue $ws.Range("D28") '168.83'
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("E29").Value = '  +0.17%  '
Set-TextValue $ws.Range("D30") '9.19'
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("E31").Value = '  -1.89%  '
Set-TextValue $ws.Range("D32") '4.96'
$ws.Range("E32").Value = '  +9.08%  '
$ws.Range("E33").Value = '  +2.37%  '
$ws.Range("E34").Value = '  -0.05%  '
Set-TextValue $ws.Range("D35") '17.95'
$ws.Range("E35").Value = '  +5.38%  '
$ws.Range("E36").Value = '  -1.43%  '
Set-TextValue $ws.Range("D37") '0.0698'
$ws.Range("E37").Value = '  +0.94%  '
$ws.Range("E38").Value = '  +1.89%  '
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("E40").Value = '  -1.02%  '
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").Value = '1.993.38'
$ws.Range("E42").Value = '  +0.04%  '
Set-TextValue $ws.Range("D43") '0.0289'
$ws.Range("E43").Value = '  +1.22%  '
Set-TextValue $ws.Range("D44") '2.24'
$ws.Range("E44").Value = '  -4.91%  '
Set-TextValue $ws.Range("D45") '10.18'
$ws.Range("E45").Value = '  +1.13%  '
Set-TextValue $ws.Range("D46") '17.61'
$ws.Range("E46").Value = '  -0.69%  '
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("B48").Value = 'BitcoinSV'
$ws.Range("C48").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue $ws.Range("D48") '76.14'
$ws.Range("E48").Value = '  +8.33%  '
$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue $ws.Range("D49") '54.94'
$ws.Range("E49").Value = '  -2.65%  '
$ws.Range("D50").Value = '2.547.30'
$ws.Range("E50").Value = '  +0.55%  '
Set-TextValue $ws.Range("D51") '1.55'
$ws.Range("E51").Value = '  +1.12%  '
